# "Generate Report for Archive" — refresh the localization-status report:
# the two handed-off files have moved from "Ready for handoff" into
# "In Translation", and the Status column(s) are re-autofit to the new
# (shorter) text on every sheet that shows a status value.

$wb = $excel.ActiveWorkbook

$newStatus = "In Translation"
# Narrower fitted width the localization tool recomputed for the
# "In Translation" label (was 17.2159881591797 for "Ready for handoff").
$fitWidth = 12.5

# --- Overview sheet: zh-cn (col E) / de-de (col F) status columns ---
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E2").Value = $newStatus
$wsOverview.Range("F2").Value = $newStatus
$wsOverview.Range("E3").Value = $newStatus
$wsOverview.Range("F3").Value = $newStatus
$wsOverview.Columns.Item(5).ColumnWidth = $fitWidth
$wsOverview.Columns.Item(6).ColumnWidth = $fitWidth

# --- zh-cn sheet: Status column (C) ---
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("C2").Value = $newStatus
$wsZhCn.Range("C3").Value = $newStatus
$wsZhCn.Columns.Item(3).ColumnWidth = $fitWidth

# --- de-de sheet: Status column (C) ---
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("C2").Value = $newStatus
$wsDeDe.Range("C3").Value = $newStatus
$wsDeDe.Columns.Item(3).ColumnWidth = $fitWidth
